# edit.ps1
#
# Commit: "Add files via upload"
#
# Appends 27 new requisition line items (worksheet rows 347-373) to the
# single sheet "Planilha1", matching the target OOXML diff:
#   - xl/sharedStrings.xml grows from 675 to 699 unique strings
#   - xl/worksheets/sheet1.xml dimension grows from A1:P346 to A1:P373
#   - 27 new <row> elements (r="347" .. r="373") are appended, reusing the
#     workbooks existing 5 cell styles (no new styles are introduced):
#       s=1 -> default text style (vertical-top)               (B,C,H,I,J,K,O,P, and blank F/G/O/P)
#       s=2 -> integer style "0;-0"                            (A,D, and populated F)
#       s=3 -> date style "m/d/yyyy" (numFmtId 14)             (E, and populated G)
#       s=4 -> currency-ish style "#,##0.00;-#,##0.00"         (L,M,N)
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 347
$lastNewRow = 373

# Insert the whole 27-row block in one shot. Excel carries the formatting
# of the row immediately above the insertion point (row 346) down onto
# every inserted row, column by column -- so A/D get the integer style,
# E/G get the date style, F gets the integer style, L/M/N get the
# currency style, and everything else gets the plain text style, all
# while reusing the 5 styles that already exist in the workbook.
$ws.Range("A" + $firstNewRow + ":P" + $lastNewRow).Insert()

# Donor range carrying the plain "style 1" (General number format,
# vertical-top alignment, no quote-prefix) used below to reset the
# OF_CDG/OF_DATA/supplier-code/supplier-name cells back to the default
# style on line items that do not have that data (matching the blank
# placeholder cells in the target rows).
$plainStyleDonor = $ws.Cells.Item(1, 2)

# Helper: write a literal string into a cell without Excels "looks like
# a number" auto-conversion mangling values such as "00000000007832"
# (which would otherwise be parsed as the number 7832, losing the
# leading zeros). Going through a self-referencing text formula and then
# collapsing it to a static value keeps the text completely literal and
# leaves the cells existing style/number-format untouched.
function Set-LiteralText($cell, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# Helper: blank out a placeholder cell, forcing its style back to the
# plain default (style 1) instead of the numeric/date style it inherited
# from the row above.
function Clear-Placeholder($cell) {
    $plainStyleDonor.Copy()
    $cell.PasteSpecial(-4122)
}

# --- row 347 (C.04.0140 / PAPEL HIGIÊNICO ( OBRA ) FARDO C/ 64 RL) ---
$ws.Cells.Item(347, 1).Value = 2212
Set-LiteralText $ws.Cells.Item(347, 2) 'IDEA INVEST. IMOBILIÁRIOS LTDA.'
Set-LiteralText $ws.Cells.Item(347, 3) 'RJ'
$ws.Cells.Item(347, 4).Value = 279
$ws.Cells.Item(347, 5).Value = 46050.7155112963
$ws.Cells.Item(347, 6).Value = 81906
$ws.Cells.Item(347, 7).Value = 46050
Set-LiteralText $ws.Cells.Item(347, 8) 'C.04.0140'
Set-LiteralText $ws.Cells.Item(347, 9) 'PAPEL HIGIÊNICO ( OBRA ) FARDO C/ 64 RL'
Set-LiteralText $ws.Cells.Item(347, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(347, 11) 'UN'
$ws.Cells.Item(347, 12).Value = 1.0
$ws.Cells.Item(347, 13).Value = 65.6
$ws.Cells.Item(347, 14).Value = 65.6
Set-LiteralText $ws.Cells.Item(347, 15) '00000000007832'
Set-LiteralText $ws.Cells.Item(347, 16) 'KLONEX'

# --- row 348 (C.04.0002 / ÁLCOOL 1 LT 92º) ---
$ws.Cells.Item(348, 1).Value = 2212
Set-LiteralText $ws.Cells.Item(348, 2) 'IDEA INVEST. IMOBILIÁRIOS LTDA.'
Set-LiteralText $ws.Cells.Item(348, 3) 'RJ'
$ws.Cells.Item(348, 4).Value = 279
$ws.Cells.Item(348, 5).Value = 46050.7155112963
$ws.Cells.Item(348, 6).Value = 81906
$ws.Cells.Item(348, 7).Value = 46050
Set-LiteralText $ws.Cells.Item(348, 8) 'C.04.0002'
Set-LiteralText $ws.Cells.Item(348, 9) 'ÁLCOOL 1 LT 92º'
Set-LiteralText $ws.Cells.Item(348, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(348, 11) 'UN'
$ws.Cells.Item(348, 12).Value = 15.0
$ws.Cells.Item(348, 13).Value = 10.85
$ws.Cells.Item(348, 14).Value = 162.75
Set-LiteralText $ws.Cells.Item(348, 15) '00000000007832'
Set-LiteralText $ws.Cells.Item(348, 16) 'KLONEX'

# --- row 349 (E.02.0016 / FITA CREPE  VERDE 48 MM X 50MM AUTOMOTIV) ---
$ws.Cells.Item(349, 1).Value = 2212
Set-LiteralText $ws.Cells.Item(349, 2) 'IDEA INVEST. IMOBILIÁRIOS LTDA.'
Set-LiteralText $ws.Cells.Item(349, 3) 'RJ'
$ws.Cells.Item(349, 4).Value = 279
$ws.Cells.Item(349, 5).Value = 46050.7155112963
$ws.Cells.Item(349, 6).Value = 81908
$ws.Cells.Item(349, 7).Value = 46050
Set-LiteralText $ws.Cells.Item(349, 8) 'E.02.0016'
Set-LiteralText $ws.Cells.Item(349, 9) 'FITA CREPE  VERDE 48 MM X 50MM AUTOMOTIVA'
Set-LiteralText $ws.Cells.Item(349, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(349, 11) 'UN'
$ws.Cells.Item(349, 12).Value = 50.0
$ws.Cells.Item(349, 13).Value = 16.0
$ws.Cells.Item(349, 14).Value = 800.0
Set-LiteralText $ws.Cells.Item(349, 15) '00000000002393'
Set-LiteralText $ws.Cells.Item(349, 16) 'SPW3'

# --- row 350 (E.04.0410 / LÂMINAS PARA REPOSIÇÃO DE ESTILETE - 18 ) ---
$ws.Cells.Item(350, 1).Value = 2212
Set-LiteralText $ws.Cells.Item(350, 2) 'IDEA INVEST. IMOBILIÁRIOS LTDA.'
Set-LiteralText $ws.Cells.Item(350, 3) 'RJ'
$ws.Cells.Item(350, 4).Value = 279
$ws.Cells.Item(350, 5).Value = 46050.7155112963
$ws.Cells.Item(350, 6).Value = 81907
$ws.Cells.Item(350, 7).Value = 46050
Set-LiteralText $ws.Cells.Item(350, 8) 'E.04.0410'
Set-LiteralText $ws.Cells.Item(350, 9) 'LÂMINAS PARA REPOSIÇÃO DE ESTILETE - 18 MM - CX C/ 10'
Set-LiteralText $ws.Cells.Item(350, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(350, 11) 'CX'
$ws.Cells.Item(350, 12).Value = 3.0
$ws.Cells.Item(350, 13).Value = 5.9
$ws.Cells.Item(350, 14).Value = 17.7
Set-LiteralText $ws.Cells.Item(350, 15) '00000000008655'
Set-LiteralText $ws.Cells.Item(350, 16) 'DMC MATERIAIS'

# --- row 351 (C.04.0010 / PANO DE CHÃO GRANDE) ---
$ws.Cells.Item(351, 1).Value = 2212
Set-LiteralText $ws.Cells.Item(351, 2) 'IDEA INVEST. IMOBILIÁRIOS LTDA.'
Set-LiteralText $ws.Cells.Item(351, 3) 'RJ'
$ws.Cells.Item(351, 4).Value = 279
$ws.Cells.Item(351, 5).Value = 46050.7155112963
Clear-Placeholder $ws.Cells.Item(351, 6)
Clear-Placeholder $ws.Cells.Item(351, 7)
Set-LiteralText $ws.Cells.Item(351, 8) 'C.04.0010'
Set-LiteralText $ws.Cells.Item(351, 9) 'PANO DE CHÃO GRANDE'
Set-LiteralText $ws.Cells.Item(351, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(351, 11) 'UN'
$ws.Cells.Item(351, 12).Value = 80.0
$ws.Cells.Item(351, 13).Value = 0.0
$ws.Cells.Item(351, 14).Value = 0.0
Clear-Placeholder $ws.Cells.Item(351, 15)
Clear-Placeholder $ws.Cells.Item(351, 16)

# --- row 352 (E.05.0002 / PROPÉ SAPATILHA DESCARTAVEL) ---
$ws.Cells.Item(352, 1).Value = 2212
Set-LiteralText $ws.Cells.Item(352, 2) 'IDEA INVEST. IMOBILIÁRIOS LTDA.'
Set-LiteralText $ws.Cells.Item(352, 3) 'RJ'
$ws.Cells.Item(352, 4).Value = 279
$ws.Cells.Item(352, 5).Value = 46050.7155112963
Clear-Placeholder $ws.Cells.Item(352, 6)
Clear-Placeholder $ws.Cells.Item(352, 7)
Set-LiteralText $ws.Cells.Item(352, 8) 'E.05.0002'
Set-LiteralText $ws.Cells.Item(352, 9) 'PROPÉ SAPATILHA DESCARTAVEL'
Set-LiteralText $ws.Cells.Item(352, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(352, 11) 'UN'
$ws.Cells.Item(352, 12).Value = 15.0
$ws.Cells.Item(352, 13).Value = 0.0
$ws.Cells.Item(352, 14).Value = 0.0
Clear-Placeholder $ws.Cells.Item(352, 15)
Clear-Placeholder $ws.Cells.Item(352, 16)

# --- row 353 (J.02.0905 / ARGAMASSA ESPECIAL  SUPERFLEX AC III - P) ---
$ws.Cells.Item(353, 1).Value = 2212
Set-LiteralText $ws.Cells.Item(353, 2) 'IDEA INVEST. IMOBILIÁRIOS LTDA.'
Set-LiteralText $ws.Cells.Item(353, 3) 'RJ'
$ws.Cells.Item(353, 4).Value = 279
$ws.Cells.Item(353, 5).Value = 46050.7155112963
Clear-Placeholder $ws.Cells.Item(353, 6)
Clear-Placeholder $ws.Cells.Item(353, 7)
Set-LiteralText $ws.Cells.Item(353, 8) 'J.02.0905'
Set-LiteralText $ws.Cells.Item(353, 9) 'ARGAMASSA ESPECIAL  SUPERFLEX AC III - PORTOKOLL 20KG BRANCA'
Set-LiteralText $ws.Cells.Item(353, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(353, 11) 'UN'
$ws.Cells.Item(353, 12).Value = 10.0
$ws.Cells.Item(353, 13).Value = 0.0
$ws.Cells.Item(353, 14).Value = 0.0
Clear-Placeholder $ws.Cells.Item(353, 15)
Clear-Placeholder $ws.Cells.Item(353, 16)

# --- row 354 (R.02.0022 / TINTA ACRILICA LATA DE 18L) ---
$ws.Cells.Item(354, 1).Value = 2212
Set-LiteralText $ws.Cells.Item(354, 2) 'IDEA INVEST. IMOBILIÁRIOS LTDA.'
Set-LiteralText $ws.Cells.Item(354, 3) 'RJ'
$ws.Cells.Item(354, 4).Value = 279
$ws.Cells.Item(354, 5).Value = 46050.7155112963
Clear-Placeholder $ws.Cells.Item(354, 6)
Clear-Placeholder $ws.Cells.Item(354, 7)
Set-LiteralText $ws.Cells.Item(354, 8) 'R.02.0022'
Set-LiteralText $ws.Cells.Item(354, 9) 'TINTA ACRILICA LATA DE 18L'
Set-LiteralText $ws.Cells.Item(354, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(354, 11) 'LAT'
$ws.Cells.Item(354, 12).Value = 1.0
$ws.Cells.Item(354, 13).Value = 0.0
$ws.Cells.Item(354, 14).Value = 0.0
Clear-Placeholder $ws.Cells.Item(354, 15)
Clear-Placeholder $ws.Cells.Item(354, 16)

# --- row 355 (S.07.0011 / HIDROFUGANTE) ---
$ws.Cells.Item(355, 1).Value = 2212
Set-LiteralText $ws.Cells.Item(355, 2) 'IDEA INVEST. IMOBILIÁRIOS LTDA.'
Set-LiteralText $ws.Cells.Item(355, 3) 'RJ'
$ws.Cells.Item(355, 4).Value = 279
$ws.Cells.Item(355, 5).Value = 46050.7155112963
Clear-Placeholder $ws.Cells.Item(355, 6)
Clear-Placeholder $ws.Cells.Item(355, 7)
Set-LiteralText $ws.Cells.Item(355, 8) 'S.07.0011'
Set-LiteralText $ws.Cells.Item(355, 9) 'HIDROFUGANTE'
Set-LiteralText $ws.Cells.Item(355, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(355, 11) 'L'
$ws.Cells.Item(355, 12).Value = 1.0
$ws.Cells.Item(355, 13).Value = 0.0
$ws.Cells.Item(355, 14).Value = 0.0
Clear-Placeholder $ws.Cells.Item(355, 15)
Clear-Placeholder $ws.Cells.Item(355, 16)

# --- row 356 (S.08.0200 / VIAPLUS 1000/TOP IMPER. BI-COMPONENTE(A+) ---
$ws.Cells.Item(356, 1).Value = 2212
Set-LiteralText $ws.Cells.Item(356, 2) 'IDEA INVEST. IMOBILIÁRIOS LTDA.'
Set-LiteralText $ws.Cells.Item(356, 3) 'RJ'
$ws.Cells.Item(356, 4).Value = 279
$ws.Cells.Item(356, 5).Value = 46050.7155112963
Clear-Placeholder $ws.Cells.Item(356, 6)
Clear-Placeholder $ws.Cells.Item(356, 7)
Set-LiteralText $ws.Cells.Item(356, 8) 'S.08.0200'
Set-LiteralText $ws.Cells.Item(356, 9) 'VIAPLUS 1000/TOP IMPER. BI-COMPONENTE(A+B)  - EMB. 18KG'
Set-LiteralText $ws.Cells.Item(356, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(356, 11) 'CX'
$ws.Cells.Item(356, 12).Value = 5.0
$ws.Cells.Item(356, 13).Value = 0.0
$ws.Cells.Item(356, 14).Value = 0.0
Clear-Placeholder $ws.Cells.Item(356, 15)
Clear-Placeholder $ws.Cells.Item(356, 16)

# --- row 357 (S.08.0203 / VIAPLUS 7000 - FIBRAS  IMPER. BI-COMPONE) ---
$ws.Cells.Item(357, 1).Value = 2212
Set-LiteralText $ws.Cells.Item(357, 2) 'IDEA INVEST. IMOBILIÁRIOS LTDA.'
Set-LiteralText $ws.Cells.Item(357, 3) 'RJ'
$ws.Cells.Item(357, 4).Value = 279
$ws.Cells.Item(357, 5).Value = 46050.7155112963
Clear-Placeholder $ws.Cells.Item(357, 6)
Clear-Placeholder $ws.Cells.Item(357, 7)
Set-LiteralText $ws.Cells.Item(357, 8) 'S.08.0203'
Set-LiteralText $ws.Cells.Item(357, 9) 'VIAPLUS 7000 - FIBRAS  IMPER. BI-COMPONENTE(A+B)  - EMB. 18KG'
Set-LiteralText $ws.Cells.Item(357, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(357, 11) 'CX'
$ws.Cells.Item(357, 12).Value = 7.0
$ws.Cells.Item(357, 13).Value = 0.0
$ws.Cells.Item(357, 14).Value = 0.0
Clear-Placeholder $ws.Cells.Item(357, 15)
Clear-Placeholder $ws.Cells.Item(357, 16)

# --- row 358 (E.04.0718 / DISCO DE CORTE PARA FERRO 4/12'') ---
$ws.Cells.Item(358, 1).Value = 2504
Set-LiteralText $ws.Cells.Item(358, 2) 'MARIA ANGÉLICA A. M. DA COSTA'
Set-LiteralText $ws.Cells.Item(358, 3) 'RJ'
$ws.Cells.Item(358, 4).Value = 47
$ws.Cells.Item(358, 5).Value = 46050.7176682523
$ws.Cells.Item(358, 6).Value = 81909
$ws.Cells.Item(358, 7).Value = 46050
Set-LiteralText $ws.Cells.Item(358, 8) 'E.04.0718'
Set-LiteralText $ws.Cells.Item(358, 9) 'DISCO DE CORTE PARA FERRO 4/12'''''
Set-LiteralText $ws.Cells.Item(358, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(358, 11) 'UN'
$ws.Cells.Item(358, 12).Value = 4.0
$ws.Cells.Item(358, 13).Value = 2.75
$ws.Cells.Item(358, 14).Value = 11.0
Set-LiteralText $ws.Cells.Item(358, 15) '00000000008655'
Set-LiteralText $ws.Cells.Item(358, 16) 'DMC MATERIAIS'

# --- row 359 (E.04.0650 / LIMA  CHATA) ---
$ws.Cells.Item(359, 1).Value = 2504
Set-LiteralText $ws.Cells.Item(359, 2) 'MARIA ANGÉLICA A. M. DA COSTA'
Set-LiteralText $ws.Cells.Item(359, 3) 'RJ'
$ws.Cells.Item(359, 4).Value = 47
$ws.Cells.Item(359, 5).Value = 46050.7176682523
Clear-Placeholder $ws.Cells.Item(359, 6)
Clear-Placeholder $ws.Cells.Item(359, 7)
Set-LiteralText $ws.Cells.Item(359, 8) 'E.04.0650'
Set-LiteralText $ws.Cells.Item(359, 9) 'LIMA  CHATA'
Set-LiteralText $ws.Cells.Item(359, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(359, 11) 'UN'
$ws.Cells.Item(359, 12).Value = 1.0
$ws.Cells.Item(359, 13).Value = 0.0
$ws.Cells.Item(359, 14).Value = 0.0
Clear-Placeholder $ws.Cells.Item(359, 15)
Clear-Placeholder $ws.Cells.Item(359, 16)

# --- row 360 (E.04.0647 / LÂMINA DE SERRA P/ ALUMÍNIO - STARRET) ---
$ws.Cells.Item(360, 1).Value = 2504
Set-LiteralText $ws.Cells.Item(360, 2) 'MARIA ANGÉLICA A. M. DA COSTA'
Set-LiteralText $ws.Cells.Item(360, 3) 'RJ'
$ws.Cells.Item(360, 4).Value = 47
$ws.Cells.Item(360, 5).Value = 46050.7176682523
Clear-Placeholder $ws.Cells.Item(360, 6)
Clear-Placeholder $ws.Cells.Item(360, 7)
Set-LiteralText $ws.Cells.Item(360, 8) 'E.04.0647'
Set-LiteralText $ws.Cells.Item(360, 9) 'LÂMINA DE SERRA P/ ALUMÍNIO - STARRET'
Set-LiteralText $ws.Cells.Item(360, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(360, 11) 'UN'
$ws.Cells.Item(360, 12).Value = 3.0
$ws.Cells.Item(360, 13).Value = 0.0
$ws.Cells.Item(360, 14).Value = 0.0
Clear-Placeholder $ws.Cells.Item(360, 15)
Clear-Placeholder $ws.Cells.Item(360, 16)

# --- row 361 (J.02.0905 / ARGAMASSA ESPECIAL  SUPERFLEX AC III - P) ---
$ws.Cells.Item(361, 1).Value = 2504
Set-LiteralText $ws.Cells.Item(361, 2) 'MARIA ANGÉLICA A. M. DA COSTA'
Set-LiteralText $ws.Cells.Item(361, 3) 'RJ'
$ws.Cells.Item(361, 4).Value = 47
$ws.Cells.Item(361, 5).Value = 46050.7176682523
Clear-Placeholder $ws.Cells.Item(361, 6)
Clear-Placeholder $ws.Cells.Item(361, 7)
Set-LiteralText $ws.Cells.Item(361, 8) 'J.02.0905'
Set-LiteralText $ws.Cells.Item(361, 9) 'ARGAMASSA ESPECIAL  SUPERFLEX AC III - PORTOKOLL 20KG BRANCA'
Set-LiteralText $ws.Cells.Item(361, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(361, 11) 'UN'
$ws.Cells.Item(361, 12).Value = 15.0
$ws.Cells.Item(361, 13).Value = 0.0
$ws.Cells.Item(361, 14).Value = 0.0
Clear-Placeholder $ws.Cells.Item(361, 15)
Clear-Placeholder $ws.Cells.Item(361, 16)

# --- row 362 (K.02.0234 / ANEL DE BORRACHA ESGOTO SERIE NORMAL 150) ---
$ws.Cells.Item(362, 1).Value = 2504
Set-LiteralText $ws.Cells.Item(362, 2) 'MARIA ANGÉLICA A. M. DA COSTA'
Set-LiteralText $ws.Cells.Item(362, 3) 'RJ'
$ws.Cells.Item(362, 4).Value = 47
$ws.Cells.Item(362, 5).Value = 46050.7176682523
Clear-Placeholder $ws.Cells.Item(362, 6)
Clear-Placeholder $ws.Cells.Item(362, 7)
Set-LiteralText $ws.Cells.Item(362, 8) 'K.02.0234'
Set-LiteralText $ws.Cells.Item(362, 9) 'ANEL DE BORRACHA ESGOTO SERIE NORMAL 150 MM'
Set-LiteralText $ws.Cells.Item(362, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(362, 11) 'UN'
$ws.Cells.Item(362, 12).Value = 21.0
$ws.Cells.Item(362, 13).Value = 0.0
$ws.Cells.Item(362, 14).Value = 0.0
Clear-Placeholder $ws.Cells.Item(362, 15)
Clear-Placeholder $ws.Cells.Item(362, 16)

# --- row 363 (K.02.0129 / TUBO PVC ESGOTO SERIE R  150MM C/ 6 M  A) ---
$ws.Cells.Item(363, 1).Value = 2504
Set-LiteralText $ws.Cells.Item(363, 2) 'MARIA ANGÉLICA A. M. DA COSTA'
Set-LiteralText $ws.Cells.Item(363, 3) 'RJ'
$ws.Cells.Item(363, 4).Value = 47
$ws.Cells.Item(363, 5).Value = 46050.7176682523
Clear-Placeholder $ws.Cells.Item(363, 6)
Clear-Placeholder $ws.Cells.Item(363, 7)
Set-LiteralText $ws.Cells.Item(363, 8) 'K.02.0129'
Set-LiteralText $ws.Cells.Item(363, 9) 'TUBO PVC ESGOTO SERIE R  150MM C/ 6 M  AMANCO'
Set-LiteralText $ws.Cells.Item(363, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(363, 11) 'UN'
$ws.Cells.Item(363, 12).Value = 1.0
$ws.Cells.Item(363, 13).Value = 0.0
$ws.Cells.Item(363, 14).Value = 0.0
Clear-Placeholder $ws.Cells.Item(363, 15)
Clear-Placeholder $ws.Cells.Item(363, 16)

# --- row 364 (K.02.3556 / JOELHO 45° PVC P/ ESG. SÉRIE R DN 150MM ) ---
$ws.Cells.Item(364, 1).Value = 2504
Set-LiteralText $ws.Cells.Item(364, 2) 'MARIA ANGÉLICA A. M. DA COSTA'
Set-LiteralText $ws.Cells.Item(364, 3) 'RJ'
$ws.Cells.Item(364, 4).Value = 47
$ws.Cells.Item(364, 5).Value = 46050.7176682523
Clear-Placeholder $ws.Cells.Item(364, 6)
Clear-Placeholder $ws.Cells.Item(364, 7)
Set-LiteralText $ws.Cells.Item(364, 8) 'K.02.3556'
Set-LiteralText $ws.Cells.Item(364, 9) 'JOELHO 45° PVC P/ ESG. SÉRIE R DN 150MM TIGRE'
Set-LiteralText $ws.Cells.Item(364, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(364, 11) 'UN'
$ws.Cells.Item(364, 12).Value = 21.0
$ws.Cells.Item(364, 13).Value = 0.0
$ws.Cells.Item(364, 14).Value = 0.0
Clear-Placeholder $ws.Cells.Item(364, 15)
Clear-Placeholder $ws.Cells.Item(364, 16)

# --- row 365 (W.06.0014 / ADESIVO ESTRUTURAL DE CONSISTENCIA FLUID) ---
$ws.Cells.Item(365, 1).Value = 2504
Set-LiteralText $ws.Cells.Item(365, 2) 'MARIA ANGÉLICA A. M. DA COSTA'
Set-LiteralText $ws.Cells.Item(365, 3) 'RJ'
$ws.Cells.Item(365, 4).Value = 47
$ws.Cells.Item(365, 5).Value = 46050.7176682523
Clear-Placeholder $ws.Cells.Item(365, 6)
Clear-Placeholder $ws.Cells.Item(365, 7)
Set-LiteralText $ws.Cells.Item(365, 8) 'W.06.0014'
Set-LiteralText $ws.Cells.Item(365, 9) 'ADESIVO ESTRUTURAL DE CONSISTENCIA FLUIDA - DENVERPOXI'
Set-LiteralText $ws.Cells.Item(365, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(365, 11) 'UN'
$ws.Cells.Item(365, 12).Value = 8.0
$ws.Cells.Item(365, 13).Value = 0.0
$ws.Cells.Item(365, 14).Value = 0.0
Clear-Placeholder $ws.Cells.Item(365, 15)
Clear-Placeholder $ws.Cells.Item(365, 16)

# --- row 366 (J.01.0016 / PEDRA BRITADA Nº 1 - SACO GRANDE 20 KG) ---
$ws.Cells.Item(366, 1).Value = 2317
Set-LiteralText $ws.Cells.Item(366, 2) 'LUIZ ALBERTO HESS BORGES'
Set-LiteralText $ws.Cells.Item(366, 3) 'SP'
$ws.Cells.Item(366, 4).Value = 67
$ws.Cells.Item(366, 5).Value = 46050.7192297917
$ws.Cells.Item(366, 6).Value = 81910
$ws.Cells.Item(366, 7).Value = 46050
Set-LiteralText $ws.Cells.Item(366, 8) 'J.01.0016'
Set-LiteralText $ws.Cells.Item(366, 9) 'PEDRA BRITADA Nº 1 - SACO GRANDE 20 KG'
Set-LiteralText $ws.Cells.Item(366, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(366, 11) 'SC'
$ws.Cells.Item(366, 12).Value = 80.0
$ws.Cells.Item(366, 13).Value = 5.8
$ws.Cells.Item(366, 14).Value = 464.0
Set-LiteralText $ws.Cells.Item(366, 15) '00000000009733'
Set-LiteralText $ws.Cells.Item(366, 16) 'CASA PEDROSO2648864-'

# --- row 367 (J.05.0001 / CIMENTO CP II - E-32 - 50 KG) ---
$ws.Cells.Item(367, 1).Value = 2317
Set-LiteralText $ws.Cells.Item(367, 2) 'LUIZ ALBERTO HESS BORGES'
Set-LiteralText $ws.Cells.Item(367, 3) 'SP'
$ws.Cells.Item(367, 4).Value = 67
$ws.Cells.Item(367, 5).Value = 46050.7192297917
$ws.Cells.Item(367, 6).Value = 81910
$ws.Cells.Item(367, 7).Value = 46050
Set-LiteralText $ws.Cells.Item(367, 8) 'J.05.0001'
Set-LiteralText $ws.Cells.Item(367, 9) 'CIMENTO CP II - E-32 - 50 KG'
Set-LiteralText $ws.Cells.Item(367, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(367, 11) 'SC'
$ws.Cells.Item(367, 12).Value = 30.0
$ws.Cells.Item(367, 13).Value = 37.95
$ws.Cells.Item(367, 14).Value = 1138.5
Set-LiteralText $ws.Cells.Item(367, 15) '00000000009733'
Set-LiteralText $ws.Cells.Item(367, 16) 'CASA PEDROSO2648864-'

# --- row 368 (C.04.0010 / PANO DE CHÃO GRANDE) ---
$ws.Cells.Item(368, 1).Value = 2317
Set-LiteralText $ws.Cells.Item(368, 2) 'LUIZ ALBERTO HESS BORGES'
Set-LiteralText $ws.Cells.Item(368, 3) 'SP'
$ws.Cells.Item(368, 4).Value = 67
$ws.Cells.Item(368, 5).Value = 46050.7192297917
Clear-Placeholder $ws.Cells.Item(368, 6)
Clear-Placeholder $ws.Cells.Item(368, 7)
Set-LiteralText $ws.Cells.Item(368, 8) 'C.04.0010'
Set-LiteralText $ws.Cells.Item(368, 9) 'PANO DE CHÃO GRANDE'
Set-LiteralText $ws.Cells.Item(368, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(368, 11) 'UN'
$ws.Cells.Item(368, 12).Value = 100.0
$ws.Cells.Item(368, 13).Value = 0.0
$ws.Cells.Item(368, 14).Value = 0.0
Clear-Placeholder $ws.Cells.Item(368, 15)
Clear-Placeholder $ws.Cells.Item(368, 16)

# --- row 369 (C.04.0151 / VEJA MULTIUSO - 500 ML) ---
$ws.Cells.Item(369, 1).Value = 2317
Set-LiteralText $ws.Cells.Item(369, 2) 'LUIZ ALBERTO HESS BORGES'
Set-LiteralText $ws.Cells.Item(369, 3) 'SP'
$ws.Cells.Item(369, 4).Value = 67
$ws.Cells.Item(369, 5).Value = 46050.7192297917
Clear-Placeholder $ws.Cells.Item(369, 6)
Clear-Placeholder $ws.Cells.Item(369, 7)
Set-LiteralText $ws.Cells.Item(369, 8) 'C.04.0151'
Set-LiteralText $ws.Cells.Item(369, 9) 'VEJA MULTIUSO - 500 ML'
Set-LiteralText $ws.Cells.Item(369, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(369, 11) 'UN'
$ws.Cells.Item(369, 12).Value = 4.0
$ws.Cells.Item(369, 13).Value = 0.0
$ws.Cells.Item(369, 14).Value = 0.0
Clear-Placeholder $ws.Cells.Item(369, 15)
Clear-Placeholder $ws.Cells.Item(369, 16)

# --- row 370 (H.11.0104 / TELA SOLDADA AÇO CA60 - Q196 - 2,45 X 6,) ---
$ws.Cells.Item(370, 1).Value = 2317
Set-LiteralText $ws.Cells.Item(370, 2) 'LUIZ ALBERTO HESS BORGES'
Set-LiteralText $ws.Cells.Item(370, 3) 'SP'
$ws.Cells.Item(370, 4).Value = 67
$ws.Cells.Item(370, 5).Value = 46050.7192297917
Clear-Placeholder $ws.Cells.Item(370, 6)
Clear-Placeholder $ws.Cells.Item(370, 7)
Set-LiteralText $ws.Cells.Item(370, 8) 'H.11.0104'
Set-LiteralText $ws.Cells.Item(370, 9) 'TELA SOLDADA AÇO CA60 - Q196 - 2,45 X 6,00 M'
Set-LiteralText $ws.Cells.Item(370, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(370, 11) 'UN'
$ws.Cells.Item(370, 12).Value = 5.0
$ws.Cells.Item(370, 13).Value = 0.0
$ws.Cells.Item(370, 14).Value = 0.0
Clear-Placeholder $ws.Cells.Item(370, 15)
Clear-Placeholder $ws.Cells.Item(370, 16)

# --- row 371 (H.11.0014 / AÇO CA25 10,0 MM - VARA) ---
$ws.Cells.Item(371, 1).Value = 2317
Set-LiteralText $ws.Cells.Item(371, 2) 'LUIZ ALBERTO HESS BORGES'
Set-LiteralText $ws.Cells.Item(371, 3) 'SP'
$ws.Cells.Item(371, 4).Value = 67
$ws.Cells.Item(371, 5).Value = 46050.7192297917
Clear-Placeholder $ws.Cells.Item(371, 6)
Clear-Placeholder $ws.Cells.Item(371, 7)
Set-LiteralText $ws.Cells.Item(371, 8) 'H.11.0014'
Set-LiteralText $ws.Cells.Item(371, 9) 'AÇO CA25 10,0 MM - VARA'
Set-LiteralText $ws.Cells.Item(371, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(371, 11) 'UN'
$ws.Cells.Item(371, 12).Value = 2.0
$ws.Cells.Item(371, 13).Value = 0.0
$ws.Cells.Item(371, 14).Value = 0.0
Clear-Placeholder $ws.Cells.Item(371, 15)
Clear-Placeholder $ws.Cells.Item(371, 16)

# --- row 372 (J.03.0015 / AREIA  - SACO GRANDE 20KG) ---
$ws.Cells.Item(372, 1).Value = 2317
Set-LiteralText $ws.Cells.Item(372, 2) 'LUIZ ALBERTO HESS BORGES'
Set-LiteralText $ws.Cells.Item(372, 3) 'SP'
$ws.Cells.Item(372, 4).Value = 67
$ws.Cells.Item(372, 5).Value = 46050.7192297917
Clear-Placeholder $ws.Cells.Item(372, 6)
Clear-Placeholder $ws.Cells.Item(372, 7)
Set-LiteralText $ws.Cells.Item(372, 8) 'J.03.0015'
Set-LiteralText $ws.Cells.Item(372, 9) 'AREIA  - SACO GRANDE 20KG'
Set-LiteralText $ws.Cells.Item(372, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(372, 11) 'SC'
$ws.Cells.Item(372, 12).Value = 300.0
$ws.Cells.Item(372, 13).Value = 0.0
$ws.Cells.Item(372, 14).Value = 0.0
Clear-Placeholder $ws.Cells.Item(372, 15)
Clear-Placeholder $ws.Cells.Item(372, 16)

# --- row 373 (W.06.1000 / ADESIVO ESTRUTURAL BASE EPÓXI MÉDIA FLUI) ---
$ws.Cells.Item(373, 1).Value = 2317
Set-LiteralText $ws.Cells.Item(373, 2) 'LUIZ ALBERTO HESS BORGES'
Set-LiteralText $ws.Cells.Item(373, 3) 'SP'
$ws.Cells.Item(373, 4).Value = 67
$ws.Cells.Item(373, 5).Value = 46050.7192297917
Clear-Placeholder $ws.Cells.Item(373, 6)
Clear-Placeholder $ws.Cells.Item(373, 7)
Set-LiteralText $ws.Cells.Item(373, 8) 'W.06.1000'
Set-LiteralText $ws.Cells.Item(373, 9) 'ADESIVO ESTRUTURAL BASE EPÓXI MÉDIA FLUIDEZ  COMPOUND ADESIVO  1KG   ( A+B)'
Set-LiteralText $ws.Cells.Item(373, 10) 'Apto'
Set-LiteralText $ws.Cells.Item(373, 11) 'UN'
$ws.Cells.Item(373, 12).Value = 5.0
$ws.Cells.Item(373, 13).Value = 0.0
$ws.Cells.Item(373, 14).Value = 0.0
Clear-Placeholder $ws.Cells.Item(373, 15)
Clear-Placeholder $ws.Cells.Item(373, 16)

$excel.CutCopyMode = 0
